# Adds pronunciation data (column C) for several rows, and fills in
# missing definitions (column B) for a handful of rows that were
# previously empty, on the "List 2" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List 2")

# --- Column C: Pronunciation -------------------------------------------------
$ws.Range("C3").Value  = "ˈvā-ˌdōs"
$ws.Range("C5").Value  = "kär-ˈtel"
$ws.Range("C7").Value  = "si-ˈzyu̇r-ə"
$ws.Range("C11").Value = "ˈcher-ē"
$ws.Range("C15").Value = "ˈpät-ˌshərd"
$ws.Range("C16").Value = "ˌpir-ə-ˈwet"
$ws.Range("C21").Value = "ˈspyüm"
$ws.Range("C22").Value = "ˈweft"
$ws.Range("C23").Value = "ˈpär-ˌflesh"
$ws.Range("C26").Value = "ˌī-dē-ˈā-shən"
$ws.Range("C28").Value = "ˌta-tər-di-ˈmāl-yən"
$ws.Range("C29").Value = "ˈslä-ləm"
$ws.Range("C31").Value = "ˈswāl"
$ws.Range("C32").Value = "i-ˈskärp-mənt"
$ws.Range("C34").Value = "ˈyü-kər"
$ws.Range("C36").Value = "ˈär-mə-ˌchu̇r"
$ws.Range("C38").Value = "krə-ˈvat"
$ws.Range("C39").Value = "ˈshōt"
$ws.Range("C40").Value = "kär-ˈtüsh"
$ws.Range("C42").Value = "ˈkre-nə-ˌlā-təd"
$ws.Range("C44").Value = "ˈär-gə-sē"
$ws.Range("C45").Value = "ˈslō"
$ws.Range("C49").Value = "ˈpül-ˌkā"
$ws.Range("C50").Value = "ˌshi-və-ˈrē"

# --- Column B: Definition (rows that were previously blank) -----------------
$ws.Range("B44").Value = "a large ship; especially : a large merchant ship"
$ws.Range("B45").Value = "the small dark globose astringent fruit of the blackthorn; also : blackthorn"
$ws.Range("B49").Value = "a Mexican alcoholic beverage made from the fermented sap of various agaves (such as Agave atrovirens)"
$ws.Range("B50").Value = "shivaree"

# Match the formatting already used by the rest of column B / column A for the
# rows whose definition cell is newly populated (or was touched) in this edit.
$ws.Range("B34").Style  = $ws.Range("A34").Style
$ws.Range("B36").Style  = $ws.Range("A36").Style
$ws.Range("B38").Style  = $ws.Range("A38").Style
$ws.Range("B39").Style  = $ws.Range("A39").Style
$ws.Range("B40").Style  = $ws.Range("A40").Style
$ws.Range("B41").Style  = $ws.Range("A41").Style
$ws.Range("B42").Style  = $ws.Range("A42").Style
$ws.Range("B45").Style  = $ws.Range("A45").Style
